$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch to manual calculation so the text values entered below don't force
# a full recalculation of the sheet (and flip unrelated formula cells such
# as B6's SUM into an error) - matches the source edit, which only patched
# the four cells below and left every other cached value untouched.
$excel.Calculation = -4135  # xlCalculationManual

# Localize the header row: "Date" -> "Fecha" and fix the casing of
# "Chai" -> "chai" in the total-sales column header.
$ws.Range("A1").Value = "Fecha"
$ws.Range("B1").Value = "Ventas totales de chai (unidades)"

# Row 6's pre-formatted and custom sales figures were mis-typed as clock
# times rather than plain numbers.
$ws.Range("D6").Value = "4:36"
$ws.Range("E6").Value = "05:17"
